# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" between "总计" and "2021-Q3" with the
#   full fund-holdings breakdown for the new quarter.
# - Update the "总计" (totals) sheet: the existing 2021-Q3 summary row moves
#   down one row, and a new summary row for 2022-Q4 is inserted above it.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift the existing data row down and insert the new one
# ---------------------------------------------------------------------

# Copy the index-cell style (A2, bold/centered) down to A3 before touching
# any values, so the new row keeps the same look as row 2.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# Row 3 becomes what used to be row 2 (2021-Q3 summary).
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.04

# Row 2 becomes the new 2022-Q4 summary.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 12
$wsTotal.Range("D2").Value = 0.44

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet (inserted right after the active
#    "2021-Q3" sheet, i.e. between "总计" and "2021-Q3").
# ---------------------------------------------------------------------

$wsQ4 = $wb.Worksheets.Add()
$wsQ4.Name = "2022-Q4"

# Header row - reuse the bold/bordered header style from "总计"!B1:D1.
$wsTotal.Range("B1:D1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Index column (A2:A13) - reuse the same bold/centered style as "总计"!A2.
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A13").PasteSpecial(-4122)

$rows = @(
    @(0,  "005416", "鹏华尊惠18个月定期开放混合A",                 "2.83", "39.19", "2.41", "0.0682", 2),
    @(1,  "009668", "鹏华安庆混合C",                               "2.79", "39.85", "2.41", "0.0672", 2),
    @(2,  "009667", "鹏华安庆混合A",                               "2.34", "39.85", "2.41", "0.0564", 2),
    @(3,  "011573", "鹏华安荣混合C",                               "1.87", "39.61", "2.42", "0.0453", 2),
    @(4,  "003166", "鹏华弘嘉灵活配置混合C",                       "0.92", "91.65", "4.63", "0.0426", 3),
    @(5,  "011572", "鹏华安荣混合A",                               "1.54", "39.61", "2.42", "0.0373", 2),
    @(6,  "009231", "鹏华安和混合C",                               "1.54", "38.20", "2.25", "0.0346", 2),
    @(7,  "003165", "鹏华弘嘉灵活配置混合A",                       "0.72", "91.65", "4.63", "0.0333", 3),
    @(8,  "009230", "鹏华安和混合A",                               "1.37", "38.20", "2.25", "0.0308", 2),
    @(9,  "010857", "宝盈祥乐一年持有期混合型证券投资基金A",       "1.00", "33.11", "1.26", "0.0126", 9),
    @(10, "005417", "鹏华尊惠18个月定期开放混合C",                 "0.36", "39.19", "2.41", "0.0087", 2),
    @(11, "010858", "宝盈祥乐一年持有期混合型证券投资基金C",       "0.06", "33.11", "1.26", "0.0008", 9)
)

$r = 2
foreach ($row in $rows) {
    $wsQ4.Cells.Item($r, 1).Value = $row[0]
    # Leading apostrophe forces text storage so codes like "005416" keep
    # their leading zeros and ratios like "2.83" stay text, not numbers.
    $wsQ4.Cells.Item($r, 2).Value = "'" + $row[1]
    $wsQ4.Cells.Item($r, 3).Value = $row[2]
    $wsQ4.Cells.Item($r, 4).Value = "'" + $row[3]
    $wsQ4.Cells.Item($r, 5).Value = "'" + $row[4]
    $wsQ4.Cells.Item($r, 6).Value = "'" + $row[5]
    $wsQ4.Cells.Item($r, 7).Value = "'" + $row[6]
    $wsQ4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
